$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 190 (shifts existing rows 190:221 down to 191:222,
# and extends the used range from A1:T221 to A1:T222).
$ws.Rows(190).Insert()

# Populate the newly inserted row 190 with the new weekly price entry.
$ws.Range("A190").Value = 4
$ws.Range("B190").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C190").Value = "Los Lagos"
$ws.Range("D190").Value = 44776
$ws.Range("E190").Value = 10
$ws.Range("F190").Value = "Fruta"
$ws.Range("G190").Value = 100108
$ws.Range("H190").Value = "Tropicales y subtropicales"
$ws.Range("I190").Value = 100108002
$ws.Range("J190").Value = "Mango"
$ws.Range("K190").Value = "Sin especificar"
$ws.Range("L190").Value = "Primera"
$ws.Range("M190").Value = 100
$ws.Range("N190").Value = 13000
$ws.Range("O190").Value = 14000
$ws.Range("P190").Value = 13500
$ws.Range("Q190").Value = "$/bandeja 4 kilos"
$ws.Range("R190").Value = "Brasil"
$ws.Range("S190").Value = 3375
$ws.Range("T190").Value = 4
